$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "BookShelves" (sheet1): append 3 new bookshelf product rows.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("BookShelves")

$ws1.Cells.Item(2, 1).Value = "Rhodes Solid Wood Bookshelf In Mahogany Finish"
$ws1.Cells.Item(2, 2).Value = "By Urban Ladder"
$ws1.Cells.Item(2, 3).Value = "₹14,755"

$ws1.Cells.Item(3, 1).Value = "Rhodes Solid Wood Bookshelf In Teak Finish"
$ws1.Cells.Item(3, 2).Value = "By Urban Ladder"
$ws1.Cells.Item(3, 3).Value = "₹14,755"

$ws1.Cells.Item(4, 1).Value = "Theodore Engineered Wood Bookshelf In Rustic Walnut Finish"
$ws1.Cells.Item(4, 2).Value = "By Urban Ladder"
$ws1.Cells.Item(4, 3).Value = "₹12,814"

# ---------------------------------------------------------------------------
# Sheet "submenuItems" (sheet2): add the sub-menu category/value pairs.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("submenuItems")

$ws2.Cells.Item(2, 1).Value = "Lounge Chairs"
$ws2.Cells.Item(2, 2).Value = "TV Units"

$ws2.Cells.Item(3, 1).Value = "Accent Chairs"
$ws2.Cells.Item(3, 2).Value = "Bookshelves"

$ws2.Cells.Item(4, 1).Value = "Recliners"
$ws2.Cells.Item(4, 2).Value = "Shoe Racks"

$ws2.Cells.Item(5, 1).Value = "Sofa Cum Bed"
$ws2.Cells.Item(5, 2).Value = "Prayer Units"

$ws2.Cells.Item(6, 1).Value = "UL Assured Picks"
$ws2.Cells.Item(6, 2).Value = "Showcases"

$ws2.Cells.Item(7, 1).Value = "Ottomans & Stools"
$ws2.Cells.Item(7, 2).Value = "Wall Shelves"

$ws2.Cells.Item(8, 1).Value = "Bean Bags"
$ws2.Cells.Item(8, 2).Value = "Entryway & Foyer"

$ws2.Cells.Item(9, 1).Value = "Benches"
$ws2.Cells.Item(9, 2).Value = "Room Divider"

$ws2.Cells.Item(10, 1).Value = "Bar Stools"
$ws2.Cells.Item(10, 2).Value = "Living Room Sets"

$ws2.Cells.Item(11, 1).Value = "Rocking Chairs"

$ws2.Cells.Item(12, 1).Value = "Gaming Chairs"

# ---------------------------------------------------------------------------
# Sheet "StudyChairs" (sheet3): blank out the brand/price of row 2
# (the leading apostrophe forces a real, empty text cell instead of
# Excel clearing the cell entirely when assigned a plain "").
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("StudyChairs")

$ws3.Cells.Item(2, 2).Value = "'"
$ws3.Cells.Item(2, 2).Style = "Normal"

$ws3.Cells.Item(2, 3).Value = "'"
$ws3.Cells.Item(2, 3).Style = "Normal"
